$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-6 from 2023-10-13 (45212) to 2023-10-22 (45221)
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45221
}
